$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 10 text values (inserted first so shared-string order matches
#     the target: "Graph Valid Tree", description, approach) ---
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Graph Valid Tree"
$ws.Range("B10").Value = "Return true if provided graph is a tree"
$ws.Range("C10").Value = "Create adjacency List of list. Use stack to iteratively DFS traverse. Keep track of visited neighbors using a set. Keep remove backpointer from the adjacency list. If numNodes == set size, return true. If current node already in the set, return false"

# --- Add hyperlinks for D8, D9 (existing rows that previously had no link)
#     and D10 (the brand-new row) ---
$ws.Hyperlinks.Add($ws.Range("D8"), "https://leetcode.com/problems/binary-tree-level-order-traversal/")
$ws.Range("D8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D9"), "https://leetcode.com/problems/validate-binary-search-tree/")
$ws.Range("D9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D10"), "https://leetcode.com/problems/graph-valid-tree/")
$ws.Range("D10").Style = "Hyperlink"

# --- Keep the "last active cell" selection in sync with the new last row ---
$null = $ws.Range("C10").Select()
